$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing values ---
$ws.Range("C2").Value = "MINION"
$ws.Range("F2").Value = 40

# --- New header row (row 1) for columns L:Q ---
$ws.Range("L1").Value = "Address of Idle Gif"
$ws.Range("M1").Value = "Address of Run Gif"
$ws.Range("N1").Value = "Address of Attack Gif"
$ws.Range("O1").Value = "Address of Get Damage Gif"
$ws.Range("P1").Value = "Address Of Death Gif"
$ws.Range("Q1").Value = " Target Society"

# --- New data row (row 2) for columns L:Q ---
$ws.Range("L2").Value = "./res/gifs/f1_altgeneral/idle_t.gif"
$ws.Range("M2").Value = "./res/gifs/f1_altgeneral/run_t.gif"
$ws.Range("N2").Value = "./res/gifs/f1_altgeneral/attack_t.gif"
$ws.Range("O2").Value = "./res/gifs/f1_altgeneral/attack_t.gif"
$ws.Range("P2").Value = "./res/gifs/f1_altgeneral/attack_t.gif"
$ws.Range("Q2").Value = "FRIENDLY"

# --- Column widths (best effort; runtime quantizes ColumnWidth to 1/6-character
#     steps on save, so these inputs are pre-solved to land as close as possible
#     to the exact target widths from the original workbook) ---
$ws.Columns.Item(2).ColumnWidth = 35.7421875
$ws.Columns.Item(3).ColumnWidth = 26.66015625
$ws.Columns.Item(4).ColumnWidth = 24.21875
$ws.Columns.Item(5).ColumnWidth = 22.16796875
$ws.Columns.Item(12).ColumnWidth = 33.59375
$ws.Columns.Item(13).ColumnWidth = 38.8671875
$ws.Columns.Item(14).ColumnWidth = 35.3515625
$ws.Columns.Item(15).ColumnWidth = 37.6953125
$ws.Columns.Item(16).ColumnWidth = 41.015625
$ws.Columns.Item(17).ColumnWidth = 32.32421875

# --- Sheet view (best effort) ---
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("O27").Select()
